$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","F","G","H","J","M","N","O")

$data = @(
    @(0.6218584934567843, 0.1255808333290247, 0.3359913593267265, 0.9814382793540233, 0.4137118716529784, 0.5649482420951557, 0.2857612890011794, 0.3459348674320069, 1.096101562761191, 1.905289578465414),
    @(0.5490857916998664, 0.1115784013191501, 0.333012261315389, 0.9802153319229632, 0.4118609893941141, 0.5678181318704461, 0.2858804709817093, 0.3229134246916274, 1.102643284451894, 1.907104797245609),
    @(0.5043429713382181, 0.1029363246645971, 0.3313160368147976, 0.9800599480283836, 0.4110542255999903, 0.5698619023312972, 0.2861075869822685, 0.308867857954219, 1.107087287862228, 1.909498277809831),
    @(0.4860960140117925, 0.0994036709620616, 0.3306583501678659, 0.9801463902052134, 0.4108082524469623, 0.5707656082216488, 0.2862388484923741, 0.3031671009468297, 1.109005859047159, 1.910795056572837),
    @(0.4830653172208201, 0.09881642355968268, 0.3305511705788007, 0.9801697896311268, 0.4107724057839235, 0.5709199481357885, 0.2862629823823895, 0.3022218901748062, 1.109330939539099, 1.911029793134162),
    @(0.504096941138954, 0.1028887260018223, 0.3313070310841226, 0.9800605074147057, 0.4110505732576897, 0.5698738031050397, 0.28610920048407, 0.3087908822343195, 1.107112726512163, 1.909514465463005),
    @(0.5967796431650072, 0.1207621876238534, 0.3349366339660236, 0.9808930154117377, 0.4130051918442916, 0.5658793420331989, 0.2857704330119617, 0.337978673794467, 1.098268532647268, 1.905649910107428),
    @(0.7780086776989492, 0.1554488912754834, 0.3431051961943297, 0.9872515285725143, 0.4194600690769335, 0.5602799452844351, 0.2863278562589144, 0.3959141696642945, 1.08431024235734, 1.908230099884491),
    @(0.9107924468736996, 0.18070102971609, 0.3497429050156029, 0.9948076915100472, 0.4258103910849087, 0.5575269442508954, 0.287483167760648, 0.4388915671849034, 1.076111625526821, 1.916337571568789),
    @(0.9711106639315403, 0.1921363140792209, 0.3528998563975705, 0.9988724293229723, 0.4290505812332839, 0.5565698907232957, 0.2881709279071245, 0.4585299942744854, 1.072827024540423, 1.921378985876657),
    @(0.9939381989250933, 0.1964588420425457, 0.3541149794204159, 1.000501885633597, 0.4303282300867579, 0.5562499248840851, 0.2884546992886641, 0.4659788605973247, 1.071647110772403, 1.92348290975707),
    @(0.9890225097093435, 0.19552825819639, 0.3538524085285957, 1.000146940008989, 0.4300508104319505, 0.5563169475251044, 0.2883925462598, 0.4643740769006257, 1.071898386123252, 1.923021122183115),
    @(0.9729889804409595, 0.1924920883174082, 0.352999432018521, 0.9990046773192205, 0.4291546781559674, 0.5565427162712382, 0.2881938063116962, 0.4591425743956208, 1.07272867234132, 1.92154817046773),
    @(0.9631661618454359, 0.190631326850081, 0.3524795156259444, 0.9983167590663413, 0.4286123723760085, 0.556686533907552, 0.2880751110129296, 0.45593970748169, 1.073245564648055, 1.920671327526719),
    @(0.9068486696723426, 0.1799526356254546, 0.3495393469682, 0.9945546770349125, 0.4256057190788596, 0.5575954309306894, 0.2874414855560161, 0.4376098818498164, 1.076335227998733, 1.916035356196772),
    @(0.8722767945277496, 0.1733880743259135, 0.3477707700078554, 0.9924074607808251, 0.4238513278851457, 0.5582286331992208, 0.2870943259619452, 0.4263873258109285, 1.078344539784574, 1.913538124850959),
    @(0.8523839711119763, 0.1696074308922846, 0.3467664727734956, 0.9912314923324743, 0.4228753159997041, 0.5586206303264873, 0.2869099158716963, 0.4199406961101388, 1.079542132897814, 1.912229146313024),
    @(0.8456472702445694, 0.1683265402287759, 0.3464286609655289, 0.990843471837934, 0.4225505305794002, 0.5587581282318865, 0.2868500996794623, 0.4177594138045961, 1.079954815004136, 1.911807815732971),
    @(0.8759578667618939, 0.1740873903391105, 0.3479576994304523, 0.9926299237276055, 0.4240346625201852, 0.5581583510645203, 0.2871297015732424, 0.427581131242377, 1.078126310350932, 1.913790775994528),
    @(0.9776987980844183, 0.1933840980277353, 0.353249439347195, 0.9993377387265667, 0.429416518034401, 0.5564752505627268, 0.2882515477863876, 0.4606788650655105, 1.072483063849297, 1.921975521958672),
    @(1.044112114198754, 0.2059502148155445, 0.3568224062251062, 1.004247606360849, 0.4332292107990554, 0.5556226689469526, 0.2891207362991253, 0.4823811883607476, 1.069167253692498, 1.928460572925019),
    @(1.008673860267038, 0.1992476888797228, 0.354905004493105, 1.001578993096771, 0.4311672392217503, 0.5560550734200547, 0.2886443880397778, 0.4707918833628568, 1.070902922015549, 1.924895367752413),
    @(0.874293706831395, 0.1737712500753616, 0.3478731497531271, 0.9925291660231466, 0.4239516753606694, 0.55819003852001, 0.2871136609868543, 0.4270413950583745, 1.078224839756039, 1.913676157728247),
    @(0.729042145910455, 0.1461052131865017, 0.3407833303470511, 0.9850250704817753, 0.4174321712280857, 0.5615556840317026, 0.2860461621919583, 0.3801677211708494, 1.087724680636711, 1.906442585080072)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $rowNum).Value = $rowVals[$j]
    }
}
